$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.778.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.66%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +13.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.613.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.813"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.285.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.639.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("E29").Value = "  -2.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.28%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0689"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0246"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000221"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.73%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +20.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0961"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "96.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.477.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
